$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.455.46"
$ws.Range("E2").Value = "  -0.65%  "

$ws.Range("D3").Value = "3.778.41"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.14"
$ws.Range("E5").Value = "  -1.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.33"
$ws.Range("E6").Value = "  +0.29%  "

$ws.Range("D7").Value = "3.777.36"
$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  -1.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("E10").Value = "  -2.34%  "

$ws.Range("E11").Value = "  +4.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.485"
$ws.Range("E12").Value = "  -1.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.93"
$ws.Range("E13").Value = "  -1.99%  "

$ws.Range("E14").Value = "  -2.83%  "

$ws.Range("D15").Value = "4.405.57"
$ws.Range("E15").Value = "  -0.44%  "

$ws.Range("D16").Value = "3.780.52"
$ws.Range("E16").Value = "  -0.42%  "

$ws.Range("D17").Value = "69.523.61"
$ws.Range("E17").Value = "  -0.66%  "

$ws.Range("E18").Value = "  -0.30%  "

$ws.Range("E19").Value = "  -3.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "510.51"
$ws.Range("E20").Value = "  +0.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.35"
$ws.Range("E21").Value = "  -2.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.38"
$ws.Range("E22").Value = "  -2.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.732"
$ws.Range("E23").Value = "  +0.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.50"
$ws.Range("E24").Value = "  -2.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.60"
$ws.Range("E25").Value = "  -0.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.84"
$ws.Range("E26").Value = "  -2.76%  "

$ws.Range("E27").Value = "  -2.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.61"
$ws.Range("E28").Value = "  -5.26%  "

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").Value = "  +1.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.98"
$ws.Range("E31").Value = "  +2.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.05"
$ws.Range("E32").Value = "  +3.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.80"
$ws.Range("E33").Value = "  -1.75%  "

$ws.Range("E34").Value = "  -0.44%  "

$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("E36").Value = "  -2.94%  "

$ws.Range("E37").Value = "  -1.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").Value = "  +4.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.340"
$ws.Range("E39").Value = "  +1.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "458.35"
$ws.Range("E40").Value = "  +8.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.07"
$ws.Range("E41").Value = "  -2.12%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.78"
$ws.Range("E42").Value = "  -2.24%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.98"
$ws.Range("E43").Value = "  +6.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.32"
$ws.Range("E44").Value = "  -3.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.59"
$ws.Range("E45").Value = "  -1.65%  "

$ws.Range("D46").Value = "2.958.56"
$ws.Range("E46").Value = "  -2.70%  "

$ws.Range("E47").Value = "  -0.50%  "

$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.03%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "139.11"
$ws.Range("E49").Value = "  +0.51%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.20"
$ws.Range("E50").Value = "  -0.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.47"
$ws.Range("E51").Value = "  -0.63%  "
